$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list (Price + Volume(1h) columns) per commit
# "Updated cryptos list on Sun May 12 20:36:48 UTC 2024 with GitHub Actions"
$ws.Range("D2").Value = "61.350.29"
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").Value = "2.925.55"
$ws.Range("E3").Value = "  +0.26%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.62"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.79"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.94"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("E10").Value = "  -1.77%  "
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "3.411.56"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "61.280.58"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("D17").Value = "2.928.10"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "431.86"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.53"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.11%  "
$ws.Range("E21").Value = "  -1.02%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "81.87"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.91"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.74"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.49%  "
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  -4.57%  "
$ws.Range("E29").Value = "  -0.50%  "
$ws.Range("E30").Value = "  -2.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.63"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("E32").Value = "  +1.51%  "
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").Value = "0.0₃0877"
$ws.Range("E34").Value = "  +3.64%  "
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "42.31"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.91%  "
$ws.Range("E42").Value = "  -2.08%  "
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").Value = "2.697.53"
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "365.83"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.14%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.64"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.53"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("E51").Value = "  -0.71%  "
